$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy formatting from existing header (A1) to match
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Updated values in B2 and D2 (slight precision changes)
$ws.Range("B2").Value = 0.02950307763024225
$ws.Range("D2").Value = 0.1292627146720758

# New Elapsed Time / CPU columns for rows 2-4
$ws.Range("G2").Value = 0.3776785511166963
$ws.Range("H2").Value = 0.968

$ws.Range("G3").Value = 0.3776785511166963
$ws.Range("H3").Value = 0.968

$ws.Range("G4").Value = 0.3776785511166963
$ws.Range("H4").Value = 0.968
